$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8191710965728389
$ws.Range("B3").Value = 0.6911928379723146
$ws.Range("B4").Value = 0.6127929849545581
$ws.Range("B5").Value = 0.8355975013180249
$ws.Range("B6").Value = 0.113646288209607
$ws.Range("B7").Value = 0.8693859472243745
$ws.Range("B8").Value = 0.8677175529710088
$ws.Range("B9").Value = 0.9281977665618217
$ws.Range("B10").Value = 0.5681212013925768

$ws.Range("A11").Value = "6_1"
$ws.Range("B11").Value = 0.959599044078598

$ws.Range("A12").Value = "6_2"
$ws.Range("B12").Value = 0.6634808853118711
